$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "queryEntityMapToClickhouseTable"

# Header row
$ws.Range("A1").Value = "test-id"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "condition"
$ws.Range("D1").Value = "domainName"
$ws.Range("E1").Value = "fields"
$ws.Range("F1").Value = "name"
$ws.Range("G1").Value = "order"
$ws.Range("H1").Value = "pageIndex"
$ws.Range("I1").Value = "pageSize"
$ws.Range("J1").Value = "timeout"
$ws.Range("K1").Value = "rspStatus"
$ws.Range("L1").Value = "rspCode"
$ws.Range("M1").Value = "rspMessage"

# Row 2
$ws.Range("A2").Value = "iot-connector-enittymaptoclickhousetable-1"
$ws.Range("B2").Value = "good request, data retrieved"
$ws.Range("F2").Value = "ClickhouseDriverSensor"

# Row 3
$ws.Range("A3").Value = "iot-connector-enittymaptoclickhousetable-2"
$ws.Range("B3").Value = "good request, data retrieved"
$ws.Range("F3").Value = "ClickhouseDriverSensor"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 2

# Row 4
$ws.Range("A4").Value = "iot-connector-enittymaptoclickhousetable-3"
$ws.Range("B4").Value = "good request, data retrieved"
$ws.Range("F4").Value = "ClickhouseDriverSensor"
$ws.Range("G4").Value = "value"

# Row 5
$ws.Range("A5").Value = "iot-connector-enittymaptoclickhousetable-4"
$ws.Range("B5").Value = "good request, data retrieved"
$ws.Range("C5").Value = "name='sensorA'"
$ws.Range("F5").Value = "ClickhouseDriverSensor"

# Row 6
$ws.Range("A6").Value = "iot-connector-enittymaptoclickhousetable-5"
$ws.Range("B6").Value = "good request, data retrieved"
$ws.Range("E6").Value = "timestamp"
$ws.Range("F6").Value = "ClickhouseDriverSensor"

# Formats: header row + data row 2 get the full-row bordered style from sheet1
$src.Range("A1:M2").Copy()
$ws.Range("A1:M2").PasteSpecial(-4122)

# Rows 3-6 only carry the style on columns A/B (matches authored sheet)
$src.Range("A3:B3").Copy()
$ws.Range("A3:B6").PasteSpecial(-4122)

$ws.Columns.Item(1).ColumnWidth = 62.5546875
$ws.Columns.Item(2).ColumnWidth = 27.6640625
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 16.5546875
$ws.Columns.Item(5).ColumnWidth = 23.109375
$ws.Columns.Item(6).ColumnWidth = 28.6640625

$ws.Activate()
$ws.Range("C5").Select()
